# Scheduled-runner market data refresh for the Carbuncle profit-tracking
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each block below refreshes
# the market-board derived columns (H:N) for a single leve/item row with the
# latest Universalis price snapshot; unrelated columns (A:G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1553.6
$ws.Range("I40").Value = 1254.909
$ws.Range("K40").Value = 1254.909
$ws.Range("M40").Value = -1079.909

$ws.Range("H74").Value = 4389.222
$ws.Range("I74").Value = 4214.7144
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4214.7144
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3278.7144
$ws.Range("N74").Value = -6872

$ws.Range("H77").Value = 4389.222
$ws.Range("I77").Value = 4214.7144
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 21073.572
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -16393.572
$ws.Range("N77").Value = -34360

$ws.Range("H113").Value = 13754.546
$ws.Range("I113").Value = 2933.3333
$ws.Range("J113").Value = 17812.5
$ws.Range("K113").Value = 2933.3333
$ws.Range("L113").Value = 17812.5
$ws.Range("M113").Value = 320.6667000000002
$ws.Range("N113").Value = -24320.5

$ws.Range("H127").Value = 55556824
$ws.Range("I127").Value = 250000540
$ws.Range("J127").Value = 1472
$ws.Range("K127").Value = 750001620
$ws.Range("L127").Value = 4416
$ws.Range("M127").Value = -749996660
$ws.Range("N127").Value = -14336

$ws.Range("H131").Value = 4391.73
$ws.Range("I131").Value = 698.0714
$ws.Range("J131").Value = 4993.0234
$ws.Range("K131").Value = 2094.2142
$ws.Range("L131").Value = 14979.0702
$ws.Range("M131").Value = 2945.7858
$ws.Range("N131").Value = -25059.0702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1269.7142
$ws.Range("I2").Value = 918.5
$ws.Range("J2").Value = 1972.1428
$ws.Range("K2").Value = 918.5
$ws.Range("L2").Value = 1972.1428
$ws.Range("M2").Value = -805.5
$ws.Range("N2").Value = -2198.1428

$ws.Range("H32").Value = 6115.109
$ws.Range("I32").Value = 4032.375
$ws.Range("K32").Value = 4032.375
$ws.Range("M32").Value = -3745.375

$ws.Range("H74").Value = 1185.4222
$ws.Range("I74").Value = 1133.7435
$ws.Range("K74").Value = 1133.7435
$ws.Range("M74").Value = -259.7435

$ws.Range("H77").Value = 1185.4222
$ws.Range("I77").Value = 1133.7435
$ws.Range("K77").Value = 5668.717500000001
$ws.Range("M77").Value = -1300.717500000001

$ws.Range("H110").Value = 40850.332
$ws.Range("I110").Value = 50585.75
$ws.Range("K110").Value = 50585.75
$ws.Range("M110").Value = -48540.75

$ws.Range("H116").Value = 1269.7142
$ws.Range("I116").Value = 918.5
$ws.Range("J116").Value = 1972.1428
$ws.Range("K116").Value = 918.5
$ws.Range("L116").Value = 1972.1428
$ws.Range("M116").Value = 1375.5
$ws.Range("N116").Value = -6560.1428

$ws.Range("H122").Value = 11113623
$ws.Range("I122").Value = 17546098
$ws.Range("J122").Value = 2984.0908
$ws.Range("K122").Value = 52638294
$ws.Range("L122").Value = 8952.2724
$ws.Range("M122").Value = -52635844
$ws.Range("N122").Value = -13852.2724

$ws.Range("H132").Value = 1605.0834
$ws.Range("I132").Value = 685.913
$ws.Range("J132").Value = 4625.2144
$ws.Range("K132").Value = 2057.739
$ws.Range("L132").Value = 13875.6432
$ws.Range("M132").Value = 472.261
$ws.Range("N132").Value = -18935.6432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1269.7142
$ws.Range("I3").Value = 918.5
$ws.Range("J3").Value = 1972.1428
$ws.Range("K3").Value = 918.5
$ws.Range("L3").Value = 1972.1428
$ws.Range("M3").Value = -804.5
$ws.Range("N3").Value = -2200.1428

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H135").Value = 44635
$ws.Range("J135").Value = 44635
$ws.Range("L135").Value = 44635
$ws.Range("N135").Value = -54775

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H140").Value = 60929
$ws.Range("J140").Value = 60929
$ws.Range("L140").Value = 60929
$ws.Range("N140").Value = -71289

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4056.3
$ws.Range("I31").Value = 965.7273
$ws.Range("J31").Value = 6812.757
$ws.Range("K31").Value = 965.7273
$ws.Range("L31").Value = 6812.757
$ws.Range("M31").Value = -670.7273
$ws.Range("N31").Value = -7402.757

$ws.Range("H34").Value = 4056.3
$ws.Range("I34").Value = 965.7273
$ws.Range("J34").Value = 6812.757
$ws.Range("K34").Value = 965.7273
$ws.Range("L34").Value = 6812.757
$ws.Range("M34").Value = -763.7273
$ws.Range("N34").Value = -7216.757

$ws.Range("H58").Value = 1615.7142
$ws.Range("I58").Value = 2102
$ws.Range("J58").Value = 1251
$ws.Range("K58").Value = 2102
$ws.Range("L58").Value = 1251
$ws.Range("M58").Value = -1899
$ws.Range("N58").Value = -1657

$ws.Range("H86").Value = 3232.2
$ws.Range("I86").Value = 3985.6667
$ws.Range("J86").Value = 2102
$ws.Range("K86").Value = 3985.6667
$ws.Range("L86").Value = 2102
$ws.Range("M86").Value = -2862.6667
$ws.Range("N86").Value = -4348

$ws.Range("H89").Value = 3232.2
$ws.Range("I89").Value = 3985.6667
$ws.Range("J89").Value = 2102
$ws.Range("K89").Value = 19928.3335
$ws.Range("L89").Value = 10510
$ws.Range("M89").Value = -14312.3335
$ws.Range("N89").Value = -21742

$ws.Range("H122").Value = 251475
$ws.Range("I122").Value = 500950
$ws.Range("K122").Value = 1502850
$ws.Range("M122").Value = -1500400

$ws.Range("H136").Value = 1615.7142
$ws.Range("I136").Value = 2102
$ws.Range("J136").Value = 1251
$ws.Range("K136").Value = 6306
$ws.Range("L136").Value = 3753
$ws.Range("M136").Value = -3756
$ws.Range("N136").Value = -8853

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1011.0101
$ws.Range("I68").Value = 766.7954999999999
$ws.Range("K68").Value = 2300.3865
$ws.Range("M68").Value = -1489.3865

$ws.Range("H71").Value = 1011.0101
$ws.Range("I71").Value = 766.7954999999999
$ws.Range("K71").Value = 6901.1595
$ws.Range("M71").Value = -2845.1595

$ws.Range("H132").Value = 1357
$ws.Range("I132").Value = 1439
$ws.Range("J132").Value = 1275
$ws.Range("K132").Value = 12951
$ws.Range("L132").Value = 11475
$ws.Range("M132").Value = -10421
$ws.Range("N132").Value = -16535

$ws.Range("H134").Value = 56385.832
$ws.Range("I134").Value = 56385.832
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 169157.496
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -164087.496
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 2000
$ws.Range("N10").Value = 2000
$ws.Range("M10").ClearContents()

$ws.Range("H102").Value = 1158.2
$ws.Range("I102").Value = 1225.25
$ws.Range("J102").Value = 890
$ws.Range("K102").Value = 1225.25
$ws.Range("L102").Value = 890
$ws.Range("M102").Value = 396.75
$ws.Range("N102").Value = -4134

$ws.Range("H122").Value = 58259.043
$ws.Range("I122").Value = 79737.34
$ws.Range("K122").Value = 239212.02
$ws.Range("M122").Value = -236762.02

$ws.Range("H132").Value = 2489.4146
$ws.Range("I132").Value = 1712.35
$ws.Range("J132").Value = 3229.476
$ws.Range("K132").Value = 5137.049999999999
$ws.Range("L132").Value = 9688.428
$ws.Range("M132").Value = -2607.049999999999
$ws.Range("N132").Value = -14748.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 27599.809
$ws.Range("I40").Value = 41076.92
$ws.Range("J40").Value = 5699.5
$ws.Range("K40").Value = 41076.92
$ws.Range("L40").Value = 5699.5
$ws.Range("M40").Value = -40940.92
$ws.Range("N40").Value = -5971.5

$ws.Range("H55").Value = 642.1
$ws.Range("I55").Value = 486.66666
$ws.Range("K55").Value = 486.66666
$ws.Range("M55").Value = -313.66666

$ws.Range("H132").Value = 4081.9575
$ws.Range("I132").Value = 4865.2856
$ws.Range("J132").Value = 3449.2693
$ws.Range("K132").Value = 14595.8568
$ws.Range("L132").Value = 10347.8079
$ws.Range("M132").Value = -12065.8568
$ws.Range("N132").Value = -15407.8079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 55052.156
$ws.Range("I126").Value = 73214.5
$ws.Range("J126").Value = 4197.6
$ws.Range("K126").Value = 219643.5
$ws.Range("L126").Value = 12592.8
$ws.Range("M126").Value = -217173.5
$ws.Range("N126").Value = -17532.8

